# Meeting 2 minutes — apply the "added to meeting notes" edit.
#
# 1) "Shipmap and FiringMap" gets split into separate runs wrapped in
#    w:proofErr (spellStart/spellEnd) markers, matching Word's automatic
#    "flag this as a possible misspelling" bookkeeping for the two
#    camel-case product names.
# 2) Three new sub-bullets are added under "Ethan:" (between it and the
#    placeholder bullet that used to follow it), and the trailing
#    bookmark + empty placeholder bullet are folded into the new last
#    sub-bullet.
#
# Because InsertXML only performs an in-place replace when the target
# Range spans one or more *whole* paragraphs, both edits are done by
# selecting the full paragraph(s) to change and replacing their content
# (including <w:pPr>) via InsertXML package XML.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Edit 1: "Shipmap and FiringMap" -> proofErr-wrapped runs
# ---------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Shipmap and FiringMap*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $bodyXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Admiral now contains only one Grid, which contains all the data necessary to generate the </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Shipmap</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>FiringMap</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> displays that correspond to it.</w:t></w:r></w:p>'
    $target.Range.InsertXML((New-PkgXml $bodyXml))
}

# ---------------------------------------------------------------------
# Edit 2: expand Ethan's bullet list with three new sub-items and fold
# the trailing bookmark/placeholder bullet into the last new sub-item.
# ---------------------------------------------------------------------

$workOn = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*whatever anyone else wants help on*") {
        $workOn = $i
        break
    }
}

if ($workOn -ne $null) {
    $ethan = $workOn + 1
    $placeholder = $workOn + 2

    $pEthan = $d.Paragraphs.Item($ethan)
    $pPlaceholder = $d.Paragraphs.Item($placeholder)

    # Sanity-check we are looking at the right paragraphs before rewriting.
    if (($pEthan.Range.Text -like "Ethan:*") -and ($pPlaceholder.Range.Text.Trim() -eq "")) {
        $startPos = $d.Paragraphs.Item($workOn).Range.Start
        $endPos = $pPlaceholder.Range.End
        $rng = $d.Range($startPos, $endPos)

        $bodyXml =
            '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
                '<w:r><w:t xml:space="preserve">Work on </w:t></w:r>' +
                '<w:r><w:t>whatever anyone else wants help on</w:t></w:r>' +
            '</w:p>' +
            '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
                '<w:r><w:t>Ethan:</w:t></w:r>' +
            '</w:p>' +
            '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
                '<w:r><w:t>Setup onclick events in setup/ index.html</w:t></w:r>' +
            '</w:p>' +
            '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
                '<w:r><w:t>Create methods for setup.html (preparing the game)</w:t></w:r>' +
            '</w:p>' +
            '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr>' +
                '<w:r><w:t xml:space="preserve">Add in </w:t></w:r>' +
                '<w:proofErr w:type="spellStart"/><w:r><w:t>sessionStorage</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
                '<w:r><w:t xml:space="preserve"> functionality between pages</w:t></w:r>' +
                '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
                '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '</w:p>'

        $rng.InsertXML((New-PkgXml $bodyXml))
    }
}
